$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.136.55'
$ws.Range("E2").Value = '  -4.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.653.52'
$ws.Range("E3").Value = '  -3.43%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.92'
$ws.Range("E5").Value = '  -3.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5068'
$ws.Range("E6").Value = '  -4.27%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2580'
$ws.Range("E8").Value = '  -3.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06415'
$ws.Range("E9").Value = '  -4.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.90'
$ws.Range("E10").Value = '  -4.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07782'
$ws.Range("E11").Value = '  +1.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.659.41'

$ws.Range("E13").Value = '  -5.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.882.21'
$ws.Range("E14").Value = '  -3.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5519'
$ws.Range("E15").Value = '  -5.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8014'
$ws.Range("E16").Value = '  -2.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.93'
$ws.Range("E17").Value = '  -6.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.158.92'
$ws.Range("E18").Value = '  -4.34%  '

$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '209.19'
$ws.Range("E20").Value = '  -6.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.409'
$ws.Range("E21").Value = '  -4.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.07'
$ws.Range("E22").Value = '  -3.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.051'
$ws.Range("E23").Value = '  +0.67%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.74'
$ws.Range("E25").Value = '  -0.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.744'
$ws.Range("E26").Value = '  +3.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1178'
$ws.Range("E27").Value = '  -2.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.974'
$ws.Range("E28").Value = '  -3.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.84'
$ws.Range("E29").Value = '  -2.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05088'
$ws.Range("E30").Value = '  -5.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.243'
$ws.Range("E31").Value = '  -3.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.338'
$ws.Range("E32").Value = '  -4.18%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.216'
$ws.Range("E33").Value = '  -6.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.565'
$ws.Range("E34").Value = '  -4.39%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.753'

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.368'
$ws.Range("E36").Value = '  -1.25%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9266'
$ws.Range("E37").Value = '  -2.50%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.164.54'
$ws.Range("E38").Value = '  +6.43%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5687'
$ws.Range("E39").Value = '  -2.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01589'
$ws.Range("E40").Value = '  -2.78%  '

$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.006'
$ws.Range("E41").Value = '  +0.16%  '

$ws.Range("B42").Value = 'mCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.555'
$ws.Range("E42").Value = '  -0.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8329'
$ws.Range("E43").Value = '  -1.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.659'
$ws.Range("E44").Value = '  -2.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.36'
$ws.Range("E45").Value = '  -0.63%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.792.37'
$ws.Range("E46").Value = '  -3.32%  '

$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("E48").Value = '  +0.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.67'
$ws.Range("E49").Value = '  -3.62%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.004'
$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.868'
$ws.Range("E51").Value = '  -2.64%  '
